$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records need to be inserted above the existing
# "Uva" (grape) rows that start at row 158. Insert 2 blank rows there,
# which pushes the old rows 158-173 down to 160-175 and grows the used
# range to A1:T175.
$ws.Rows("158:159").Insert()

# Row 158: new "Red Globe" record dated 45021
$ws.Range("A158").Value = 1
$ws.Range("B158").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C158").Value = "Arica y Parinacota"
$ws.Range("D158").Value = 45021
$ws.Range("E158").Value = 15
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100109
$ws.Range("H158").Value = "Uva"
$ws.Range("I158").Value = 100109001
$ws.Range("J158").Value = "Uva"
$ws.Range("K158").Value = "Red Globe"
$ws.Range("L158").Value = "Segunda"
$ws.Range("M158").Value = 300
$ws.Range("N158").Value = 16000
$ws.Range("O158").Value = 17000
$ws.Range("P158").Value = 16500
$ws.Range("Q158").Value = "$/bandeja 18 kilos"
$ws.Range("R158").Value = "Región de Coquimbo"
$ws.Range("S158").Value = 917
$ws.Range("T158").Value = 18

# Row 159: new "Thompson seedless" record dated 45021
$ws.Range("A159").Value = 1
$ws.Range("B159").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C159").Value = "Arica y Parinacota"
$ws.Range("D159").Value = 45021
$ws.Range("E159").Value = 15
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100109
$ws.Range("H159").Value = "Uva"
$ws.Range("I159").Value = 100109001
$ws.Range("J159").Value = "Uva"
$ws.Range("K159").Value = "Thompson seedless"
$ws.Range("L159").Value = "Tercera"
$ws.Range("M159").Value = 300
$ws.Range("N159").Value = 15000
$ws.Range("O159").Value = 16000
$ws.Range("P159").Value = 15500
$ws.Range("Q159").Value = "$/bandeja 18 kilos"
$ws.Range("R159").Value = "Región de Coquimbo"
$ws.Range("S159").Value = 861
$ws.Range("T159").Value = 18
